# Entsoe Consumption_Actual_Historical.xlsx - retraining the model for Horeco
# Shift the fetched window forward by 2 days (23-24 Jan 2026 -> 25-26 Jan 2026)
# and replace the "Actual Consumption (MW)" readings with the freshly retrained values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Actual Consumption (MW)" values for rows 2..193 (quarter-hours of the 2 fetched days)
$newConsumption = @(5845,5804,5768,5679,5690,5628,5634,5606,5629,5609,5620,5657,5719,5676,5730,5711,5728,5718,5765,5747,5820,5895,5873,6019,6023,6098,6089,6105,6215,6271,6322,6406,6378,6433,6468,6438,6496,6424,6396,6389,6344,6356,6285,6327,6386,6516,6543,6512,6543,6596,6652,6657,6571,6637,6619,6612,6501,6566,6603,6659,6736,6776,6951,6987,6992,7150,7315,7399,7446,7453,7491,7457,7486,7440,7415,7452,7357,7361,7322,7312,7199,7071,7054,6972,6739,6581,6432,6356,6243,6186,6041,6005,5874,5844,5831,5746,5723,5712,5692,5629,5598,5619,5596,5616,5596,5593,5590,5618,5651,5650,5686,5725,5821,5781,5892,5978,6299,6494,6678,6973,7330,7621,7735,7886,8117,8228,8411,8485,8566,8628,8673,8711,8627,8617,8558,8571,8556,8525,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

$firstRow = 2
$lastRow = 193

for ($r = $firstRow; $r -le $lastRow; $r++) {
    # --- Column A (Timestamp): roll the date window forward by 2 days, keep time-of-day ---
    $oldSerial = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 1).Value = $oldSerial + 2

    # --- Column B (Actual Consumption (MW)): new retrained values ---
    $idx = $r - $firstRow
    $ws.Cells.Item($r, 2).Value = $newConsumption[$idx]

    # --- Column D (Lookup): "dd.MM.yyyy" of the (now shifted) date + the Quarter index ---
    $newDate = $ws.Cells.Item($r, 1).Value()
    $quarter = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 4).Value = $newDate.ToString("dd.MM.yyyy") + $quarter
}

"Updated rows $firstRow..$lastRow (Timestamp +2 days, Consumption refreshed, Lookup recomputed)"
